{"js": "// Split the single \"#7cc867#fb5b89#f9cd59\" highlights paragraph into three\n// separate paragraphs, each carrying its own highlight-count summary:\n//   #7cc867: 29\n//   #fb5b89: 29\n//   #f9cd59: 19\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"#7cc867#fb5b89#f9cd59\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Insert the two new trailing paragraphs right after the target first\n  // (so the target paragraph's own position/reference stays valid), then\n  // overwrite the target paragraph's text in place for the first value.\n  target.insertParagraph(\"#f9cd59: 19\", \"After\");\n  target.insertParagraph(\"#fb5b89: 29\", \"After\");\n  target.insertText(\"#7cc867: 29\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Split the single \"#7cc867#fb5b89#f9cd59\" highlights paragraph into three\n# separate paragraphs, each carrying its own highlight-count summary:\n#   #7cc867: 29\n#   #fb5b89: 29\n#   #f9cd59: 19\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*#7cc867#fb5b89#f9cd59*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $idx = $target.Index\n\n    $r1 = $d.Paragraphs.Item($idx).Range\n    $r1.Text = \"#7cc867: 29\"\n    $r1.InsertParagraphAfter()\n\n    $r2 = $d.Paragraphs.Item($idx + 1).Range\n    $r2.Text = \"#fb5b89: 29\"\n    $r2.InsertParagraphAfter()\n\n    $r3 = $d.Paragraphs.Item($idx + 2).Range\n    $r3.Text = \"#f9cd59: 19\"\n}\n"}
